$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.421.82'
$ws.Range("E2").Value = '  -0.52%  '

$ws.Range("D3").Value = '1.924.21'
$ws.Range("E3").Value = '  +3.70%  '

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = '  +0.33%  '

$ws.Range("D5").Value = "'239.62"
$ws.Range("E5").Value = '  +2.47%  '

$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = '  +0.28%  '

$ws.Range("D7").Value = "'0.4729"
$ws.Range("E7").Value = '  -0.15%  '

$ws.Range("D8").Value = "'0.2842"
$ws.Range("E8").Value = '  +3.29%  '

$ws.Range("D9").Value = "'0.06546"
$ws.Range("E9").Value = '  +3.53%  '

$ws.Range("D10").Value = "'19.01"
$ws.Range("E10").Value = '  +7.34%  '

$ws.Range("D11").Value = "'103.74"
$ws.Range("E11").Value = '  +22.64%  '

$ws.Range("D12").Value = '1.915.43'
$ws.Range("E12").Value = '  +3.08%  '

$ws.Range("D13").Value = "'0.07579"
$ws.Range("E13").Value = '  +1.75%  '

$ws.Range("D14").Value = "'5.099"
$ws.Range("E14").Value = '  +1.95%  '

$ws.Range("D15").Value = "'0.6484"
$ws.Range("E15").Value = '  +3.56%  '

$ws.Range("D16").Value = "'292.62"
$ws.Range("E16").Value = '  +19.12%  '

$ws.Range("D17").Value = '30.423.25'
$ws.Range("E17").Value = '  -0.37%  '

$ws.Range("E18").Value = '  +0.28%  '

$ws.Range("D19").Value = "'12.93"
$ws.Range("E19").Value = '  +1.72%  '

$ws.Range("D20").Value = '2.183.24'
$ws.Range("E20").Value = '  +4.68%  '

$ws.Range("D21").Value = "'0.000007485"
$ws.Range("E21").Value = '  +2.00%  '

$ws.Range("D22").Value = "'1.005"
$ws.Range("E22").Value = '  +0.55%  '

$ws.Range("D23").Value = "'5.180"
$ws.Range("E23").Value = '  +4.82%  '

$ws.Range("D24").Value = "'6.266"
$ws.Range("E24").Value = '  +5.38%  '

$ws.Range("D25").Value = "'9.227"
$ws.Range("E25").Value = '  +0.92%  '

$ws.Range("D26").Value = "'165.53"
$ws.Range("E26").Value = '  +1.62%  '

$ws.Range("D27").Value = "'19.41"
$ws.Range("E27").Value = '  +7.66%  '

$ws.Range("D28").Value = "'2.024"
$ws.Range("E28").Value = '  +7.65%  '

$ws.Range("D29").Value = "'0.1118"
$ws.Range("E29").Value = '  +9.77%  '

$ws.Range("D30").Value = "'1.358"
$ws.Range("E30").Value = '  -0.02%  '

$ws.Range("D31").Value = "'4.091"
$ws.Range("E31").Value = '  +1.97%  '

$ws.Range("D32").Value = "'3.908"
$ws.Range("E32").Value = '  +1.84%  '

$ws.Range("D33").Value = "'0.05005"
$ws.Range("E33").Value = '  +3.39%  '

$ws.Range("D34").Value = "'0.7341"
$ws.Range("E34").Value = '  +4.48%  '

$ws.Range("D35").Value = "'1.143"
$ws.Range("E35").Value = '  +0.47%  '

$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = "'2.719"
$ws.Range("E36").Value = '  +1.04%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = "'0.01942"
$ws.Range("E37").Value = '  +2.37%  '

$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = "'2.696"
$ws.Range("E38").Value = '  +0.69%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = "'2.007"
$ws.Range("E39").Value = '  +0.33%  '

$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = "'0.8702"
$ws.Range("E40").Value = '  -0.71%  '

$ws.Range("B41").Value = 'Quant'
$ws.Range("C41").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D41").Value = "'106.96"
$ws.Range("E41").Value = '  +0.07%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = "'5.820"
$ws.Range("E42").Value = '  +4.85%  '

$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").Value = "'1.002"
$ws.Range("E43").Value = '  +0.28%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Value = "'68.73"
$ws.Range("E44").Value = '  +9.43%  '

$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = "'0.4111"
$ws.Range("E45").Value = '  +1.23%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = "'7.204"
$ws.Range("E46").Value = '  +0.10%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = "'9.236"
$ws.Range("E47").Value = '  +8.29%  '

$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").Value = "'0.1201"
$ws.Range("E48").Value = '  -0.91%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = "'34.45"
$ws.Range("E49").Value = '  +2.55%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = "'0.05632"
$ws.Range("E50").Value = '  +1.75%  '

$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = "'0.3799"
$ws.Range("E51").Value = '  +3.20%  '
